$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41..91 down to 42..92
$ws.Range("A41:R41").EntireRow.Insert()

# Populate the newly inserted row 41 with the new data record
$ws.Cells.Item(41, 1).Value = 11
$ws.Cells.Item(41, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(41, 3).Value = "Bíobío"
$ws.Cells.Item(41, 4).Value = 45175
$ws.Cells.Item(41, 4).NumberFormat = $ws.Cells.Item(42, 4).NumberFormat
$ws.Cells.Item(41, 5).Value = 8
$ws.Cells.Item(41, 6).Value = 100112031
$ws.Cells.Item(41, 7).Value = "Poroto verde"
$ws.Cells.Item(41, 8).Value = "Magnum"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 50
$ws.Cells.Item(41, 11).Value = 17000
$ws.Cells.Item(41, 12).Value = 17000
$ws.Cells.Item(41, 13).Value = 17000
$ws.Cells.Item(41, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(41, 15).Value = "Perú"
$ws.Cells.Item(41, 16).Value = 680
$ws.Cells.Item(41, 17).Value = 25
$ws.Cells.Item(41, 18).Value = "Hortaliza"
